$wb = $excel.ActiveWorkbook

# ---- Sheet 1: LP1912 ----
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 1).Value = 'Última actualización: 14:23:38'
$ws.Cells.Item(3, 1).Value = 'Total filas: 291'

# swap pairs (tie-break reorderings)
$ws.Cells.Item(14, 1).Value = '04:44:55'
$ws.Cells.Item(14, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(14, 4).Value = 2
$ws.Cells.Item(15, 1).Value = '03:46:12'
$ws.Cells.Item(15, 3).Value = '215A_EL PATO'
$ws.Cells.Item(15, 4).Value = 60
$ws.Cells.Item(67, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(68, 3).Value = '14_ABASTO'
$ws.Cells.Item(114, 3).Value = '15_ABASTO'
$ws.Cells.Item(115, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(123, 1).Value = '07:51:22'
$ws.Cells.Item(123, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(123, 4).Value = 100
$ws.Cells.Item(124, 1).Value = '08:22:12'
$ws.Cells.Item(124, 3).Value = '10_OLMOS'
$ws.Cells.Item(124, 4).Value = 69
$ws.Cells.Item(151, 1).Value = '09:57:03'
$ws.Cells.Item(151, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(151, 4).Value = 37
$ws.Cells.Item(152, 1).Value = '08:41:14'
$ws.Cells.Item(152, 3).Value = '10_OLMOS'
$ws.Cells.Item(152, 4).Value = 113

# tail rows rewrite (insert of newly-scraped rows + shift)
$ws.Cells.Item(266, 1).Value = '14:23:38'
$ws.Cells.Item(266, 2).Value = '14:33'
$ws.Cells.Item(266, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(266, 4).Value = 10
$ws.Cells.Item(266, 5).Value = 'LP1912'
$ws.Cells.Item(267, 1).Value = '12:58:39'
$ws.Cells.Item(267, 2).Value = '14:33'
$ws.Cells.Item(267, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(267, 4).Value = 95
$ws.Cells.Item(267, 5).Value = 'LP1912'
$ws.Cells.Item(268, 1).Value = '13:35:25'
$ws.Cells.Item(268, 2).Value = '14:34'
$ws.Cells.Item(268, 3).Value = '10_OLMOS'
$ws.Cells.Item(268, 4).Value = 59
$ws.Cells.Item(268, 5).Value = 'LP1912'
$ws.Cells.Item(269, 1).Value = '12:46:01'
$ws.Cells.Item(269, 2).Value = '14:34'
$ws.Cells.Item(269, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(269, 4).Value = 108
$ws.Cells.Item(269, 5).Value = 'LP1912'
$ws.Cells.Item(270, 1).Value = '12:41:18'
$ws.Cells.Item(270, 2).Value = '14:37'
$ws.Cells.Item(270, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(270, 4).Value = 116
$ws.Cells.Item(270, 5).Value = 'LP1912'
$ws.Cells.Item(271, 1).Value = '13:35:25'
$ws.Cells.Item(271, 2).Value = '14:38'
$ws.Cells.Item(271, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(271, 4).Value = 63
$ws.Cells.Item(271, 5).Value = 'LP1912'
$ws.Cells.Item(272, 1).Value = '12:41:18'
$ws.Cells.Item(272, 2).Value = '14:40'
$ws.Cells.Item(272, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(272, 4).Value = 119
$ws.Cells.Item(272, 5).Value = 'LP1912'
$ws.Cells.Item(273, 1).Value = '12:58:39'
$ws.Cells.Item(273, 2).Value = '14:44'
$ws.Cells.Item(273, 3).Value = '215B_EL PATO'
$ws.Cells.Item(273, 4).Value = 106
$ws.Cells.Item(273, 5).Value = 'LP1912'
$ws.Cells.Item(274, 1).Value = '12:46:01'
$ws.Cells.Item(274, 2).Value = '14:45'
$ws.Cells.Item(274, 3).Value = '215B_EL PATO'
$ws.Cells.Item(274, 4).Value = 119
$ws.Cells.Item(274, 5).Value = 'LP1912'
$ws.Cells.Item(275, 1).Value = '14:23:38'
$ws.Cells.Item(275, 2).Value = '14:49'
$ws.Cells.Item(275, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(275, 4).Value = 26
$ws.Cells.Item(275, 5).Value = 'LP1912'
$ws.Cells.Item(276, 1).Value = '12:58:39'
$ws.Cells.Item(276, 2).Value = '14:53'
$ws.Cells.Item(276, 3).Value = '14_ABASTO'
$ws.Cells.Item(276, 4).Value = 115
$ws.Cells.Item(276, 5).Value = 'LP1912'
$ws.Cells.Item(277, 1).Value = '12:58:39'
$ws.Cells.Item(277, 2).Value = '14:53'
$ws.Cells.Item(277, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(277, 4).Value = 115
$ws.Cells.Item(277, 5).Value = 'LP1912'
$ws.Cells.Item(278, 1).Value = '13:35:25'
$ws.Cells.Item(278, 2).Value = '14:56'
$ws.Cells.Item(278, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(278, 4).Value = 81
$ws.Cells.Item(278, 5).Value = 'LP1912'
$ws.Cells.Item(279, 1).Value = '13:35:25'
$ws.Cells.Item(279, 2).Value = '15:01'
$ws.Cells.Item(279, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(279, 4).Value = 86
$ws.Cells.Item(279, 5).Value = 'LP1912'
$ws.Cells.Item(280, 1).Value = '13:54:15'
$ws.Cells.Item(280, 2).Value = '15:02'
$ws.Cells.Item(280, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(280, 4).Value = 68
$ws.Cells.Item(280, 5).Value = 'LP1912'
$ws.Cells.Item(281, 1).Value = '14:23:38'
$ws.Cells.Item(281, 2).Value = '15:03'
$ws.Cells.Item(281, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(281, 4).Value = 40
$ws.Cells.Item(281, 5).Value = 'LP1912'
$ws.Cells.Item(282, 1).Value = '13:35:25'
$ws.Cells.Item(282, 2).Value = '15:04'
$ws.Cells.Item(282, 3).Value = '14_ABASTO'
$ws.Cells.Item(282, 4).Value = 89
$ws.Cells.Item(282, 5).Value = 'LP1912'
$ws.Cells.Item(283, 1).Value = '14:23:38'
$ws.Cells.Item(283, 2).Value = '15:04'
$ws.Cells.Item(283, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(283, 4).Value = 41
$ws.Cells.Item(283, 5).Value = 'LP1912'
$ws.Cells.Item(284, 1).Value = '13:54:15'
$ws.Cells.Item(284, 2).Value = '15:05'
$ws.Cells.Item(284, 3).Value = '14_ABASTO'
$ws.Cells.Item(284, 4).Value = 71
$ws.Cells.Item(284, 5).Value = 'LP1912'
$ws.Cells.Item(285, 1).Value = '13:35:25'
$ws.Cells.Item(285, 2).Value = '15:17'
$ws.Cells.Item(285, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(285, 4).Value = 102
$ws.Cells.Item(285, 5).Value = 'LP1912'
$ws.Cells.Item(286, 1).Value = '13:35:25'
$ws.Cells.Item(286, 2).Value = '15:24'
$ws.Cells.Item(286, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(286, 4).Value = 109
$ws.Cells.Item(286, 5).Value = 'LP1912'
$ws.Cells.Item(287, 1).Value = '14:23:38'
$ws.Cells.Item(287, 2).Value = '15:24'
$ws.Cells.Item(287, 3).Value = '215C_EL PATO'
$ws.Cells.Item(287, 4).Value = 61
$ws.Cells.Item(287, 5).Value = 'LP1912'
$ws.Cells.Item(288, 1).Value = '13:35:25'
$ws.Cells.Item(288, 2).Value = '15:25'
$ws.Cells.Item(288, 3).Value = '215C_EL PATO'
$ws.Cells.Item(288, 4).Value = 110
$ws.Cells.Item(288, 5).Value = 'LP1912'
$ws.Cells.Item(289, 1).Value = '13:54:15'
$ws.Cells.Item(289, 2).Value = '15:25'
$ws.Cells.Item(289, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(289, 4).Value = 91
$ws.Cells.Item(289, 5).Value = 'LP1912'
$ws.Cells.Item(290, 1).Value = '13:54:15'
$ws.Cells.Item(290, 2).Value = '15:25'
$ws.Cells.Item(290, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(290, 4).Value = 91
$ws.Cells.Item(290, 5).Value = 'LP1912'
$ws.Cells.Item(291, 1).Value = '13:54:15'
$ws.Cells.Item(291, 2).Value = '15:36'
$ws.Cells.Item(291, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(291, 4).Value = 102
$ws.Cells.Item(291, 5).Value = 'LP1912'
$ws.Cells.Item(292, 1).Value = '14:23:38'
$ws.Cells.Item(292, 2).Value = '15:44'
$ws.Cells.Item(292, 3).Value = '14_ABASTO'
$ws.Cells.Item(292, 4).Value = 81
$ws.Cells.Item(292, 5).Value = 'LP1912'
$ws.Cells.Item(293, 1).Value = '14:23:38'
$ws.Cells.Item(293, 2).Value = '15:45'
$ws.Cells.Item(293, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(293, 4).Value = 82
$ws.Cells.Item(293, 5).Value = 'LP1912'
$ws.Cells.Item(294, 1).Value = '14:23:38'
$ws.Cells.Item(294, 2).Value = '15:55'
$ws.Cells.Item(294, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(294, 4).Value = 92
$ws.Cells.Item(294, 5).Value = 'LP1912'
$ws.Cells.Item(295, 1).Value = '14:23:38'
$ws.Cells.Item(295, 2).Value = '16:01'
$ws.Cells.Item(295, 3).Value = '15_ABASTO'
$ws.Cells.Item(295, 4).Value = 98
$ws.Cells.Item(295, 5).Value = 'LP1912'
$ws.Cells.Item(296, 1).Value = '14:23:38'
$ws.Cells.Item(296, 2).Value = '16:20'
$ws.Cells.Item(296, 3).Value = '10_OLMOS'
$ws.Cells.Item(296, 4).Value = 117
$ws.Cells.Item(296, 5).Value = 'LP1912'

# ---- Sheet 2: LP1912-215 ----
$ws = $wb.Worksheets.Item(2)

$ws.Cells.Item(2, 1).Value = 'Última actualización: 14:23:38'
$ws.Cells.Item(3, 1).Value = 'Total filas: 76'

# swap pairs (tie-break reorderings)
$ws.Cells.Item(52, 1).Value = '11:45:06'
$ws.Cells.Item(52, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(52, 4).Value = 0
$ws.Cells.Item(53, 1).Value = '09:57:03'
$ws.Cells.Item(53, 3).Value = '215B_EL PATO'
$ws.Cells.Item(53, 4).Value = 108

# tail rows rewrite (insert of newly-scraped rows + shift)
$ws.Cells.Item(79, 1).Value = '14:23:38'
$ws.Cells.Item(79, 2).Value = '15:24'
$ws.Cells.Item(79, 3).Value = '215C_EL PATO'
$ws.Cells.Item(79, 4).Value = 61
$ws.Cells.Item(79, 5).Value = 'LP1912'
$ws.Cells.Item(80, 1).Value = '13:35:25'
$ws.Cells.Item(80, 2).Value = '15:25'
$ws.Cells.Item(80, 3).Value = '215C_EL PATO'
$ws.Cells.Item(80, 4).Value = 110
$ws.Cells.Item(80, 5).Value = 'LP1912'
$ws.Cells.Item(81, 1).Value = '14:23:38'
$ws.Cells.Item(81, 2).Value = '15:45'
$ws.Cells.Item(81, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(81, 4).Value = 82
$ws.Cells.Item(81, 5).Value = 'LP1912'

# ---- Sheet 3: 6203-6173 ----
$ws = $wb.Worksheets.Item(3)

$ws.Cells.Item(2, 1).Value = 'Última actualización: 14:23:38'
$ws.Cells.Item(3, 1).Value = 'Total filas: 291'

# swap pairs (tie-break reorderings)
$ws.Cells.Item(14, 1).Value = '04:44:55'
$ws.Cells.Item(14, 3).Value = '215_EL PELIGRO'
$ws.Cells.Item(14, 4).Value = 2
$ws.Cells.Item(15, 1).Value = '03:46:12'
$ws.Cells.Item(15, 3).Value = '215A_EL PATO'
$ws.Cells.Item(15, 4).Value = 60
$ws.Cells.Item(67, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(68, 3).Value = '14_ABASTO'
$ws.Cells.Item(114, 3).Value = '15_ABASTO'
$ws.Cells.Item(115, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(123, 1).Value = '07:51:22'
$ws.Cells.Item(123, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(123, 4).Value = 100
$ws.Cells.Item(124, 1).Value = '08:22:12'
$ws.Cells.Item(124, 3).Value = '10_OLMOS'
$ws.Cells.Item(124, 4).Value = 69
$ws.Cells.Item(151, 1).Value = '09:57:03'
$ws.Cells.Item(151, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(151, 4).Value = 37
$ws.Cells.Item(152, 1).Value = '08:41:14'
$ws.Cells.Item(152, 3).Value = '10_OLMOS'
$ws.Cells.Item(152, 4).Value = 113

# tail rows rewrite (insert of newly-scraped rows + shift)
$ws.Cells.Item(266, 1).Value = '14:23:38'
$ws.Cells.Item(266, 2).Value = '14:33'
$ws.Cells.Item(266, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(266, 4).Value = 10
$ws.Cells.Item(266, 5).Value = 'LP1912'
$ws.Cells.Item(267, 1).Value = '12:58:39'
$ws.Cells.Item(267, 2).Value = '14:33'
$ws.Cells.Item(267, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(267, 4).Value = 95
$ws.Cells.Item(267, 5).Value = 'LP1912'
$ws.Cells.Item(268, 1).Value = '13:35:25'
$ws.Cells.Item(268, 2).Value = '14:34'
$ws.Cells.Item(268, 3).Value = '10_OLMOS'
$ws.Cells.Item(268, 4).Value = 59
$ws.Cells.Item(268, 5).Value = 'LP1912'
$ws.Cells.Item(269, 1).Value = '12:46:01'
$ws.Cells.Item(269, 2).Value = '14:34'
$ws.Cells.Item(269, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(269, 4).Value = 108
$ws.Cells.Item(269, 5).Value = 'LP1912'
$ws.Cells.Item(270, 1).Value = '12:41:18'
$ws.Cells.Item(270, 2).Value = '14:37'
$ws.Cells.Item(270, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(270, 4).Value = 116
$ws.Cells.Item(270, 5).Value = 'LP1912'
$ws.Cells.Item(271, 1).Value = '13:35:25'
$ws.Cells.Item(271, 2).Value = '14:38'
$ws.Cells.Item(271, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(271, 4).Value = 63
$ws.Cells.Item(271, 5).Value = 'LP1912'
$ws.Cells.Item(272, 1).Value = '12:41:18'
$ws.Cells.Item(272, 2).Value = '14:40'
$ws.Cells.Item(272, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(272, 4).Value = 119
$ws.Cells.Item(272, 5).Value = 'LP1912'
$ws.Cells.Item(273, 1).Value = '12:58:39'
$ws.Cells.Item(273, 2).Value = '14:44'
$ws.Cells.Item(273, 3).Value = '215B_EL PATO'
$ws.Cells.Item(273, 4).Value = 106
$ws.Cells.Item(273, 5).Value = 'LP1912'
$ws.Cells.Item(274, 1).Value = '12:46:01'
$ws.Cells.Item(274, 2).Value = '14:45'
$ws.Cells.Item(274, 3).Value = '215B_EL PATO'
$ws.Cells.Item(274, 4).Value = 119
$ws.Cells.Item(274, 5).Value = 'LP1912'
$ws.Cells.Item(275, 1).Value = '14:23:38'
$ws.Cells.Item(275, 2).Value = '14:49'
$ws.Cells.Item(275, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(275, 4).Value = 26
$ws.Cells.Item(275, 5).Value = 'LP1912'
$ws.Cells.Item(276, 1).Value = '12:58:39'
$ws.Cells.Item(276, 2).Value = '14:53'
$ws.Cells.Item(276, 3).Value = '14_ABASTO'
$ws.Cells.Item(276, 4).Value = 115
$ws.Cells.Item(276, 5).Value = 'LP1912'
$ws.Cells.Item(277, 1).Value = '12:58:39'
$ws.Cells.Item(277, 2).Value = '14:53'
$ws.Cells.Item(277, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(277, 4).Value = 115
$ws.Cells.Item(277, 5).Value = 'LP1912'
$ws.Cells.Item(278, 1).Value = '13:35:25'
$ws.Cells.Item(278, 2).Value = '14:56'
$ws.Cells.Item(278, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(278, 4).Value = 81
$ws.Cells.Item(278, 5).Value = 'LP1912'
$ws.Cells.Item(279, 1).Value = '13:35:25'
$ws.Cells.Item(279, 2).Value = '15:01'
$ws.Cells.Item(279, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(279, 4).Value = 86
$ws.Cells.Item(279, 5).Value = 'LP1912'
$ws.Cells.Item(280, 1).Value = '13:54:15'
$ws.Cells.Item(280, 2).Value = '15:02'
$ws.Cells.Item(280, 3).Value = '215A_LA PLATA'
$ws.Cells.Item(280, 4).Value = 68
$ws.Cells.Item(280, 5).Value = 'LP1912'
$ws.Cells.Item(281, 1).Value = '14:23:38'
$ws.Cells.Item(281, 2).Value = '15:03'
$ws.Cells.Item(281, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(281, 4).Value = 40
$ws.Cells.Item(281, 5).Value = 'LP1912'
$ws.Cells.Item(282, 1).Value = '13:35:25'
$ws.Cells.Item(282, 2).Value = '15:04'
$ws.Cells.Item(282, 3).Value = '14_ABASTO'
$ws.Cells.Item(282, 4).Value = 89
$ws.Cells.Item(282, 5).Value = 'LP1912'
$ws.Cells.Item(283, 1).Value = '14:23:38'
$ws.Cells.Item(283, 2).Value = '15:04'
$ws.Cells.Item(283, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(283, 4).Value = 41
$ws.Cells.Item(283, 5).Value = 'LP1912'
$ws.Cells.Item(284, 1).Value = '13:54:15'
$ws.Cells.Item(284, 2).Value = '15:05'
$ws.Cells.Item(284, 3).Value = '14_ABASTO'
$ws.Cells.Item(284, 4).Value = 71
$ws.Cells.Item(284, 5).Value = 'LP1912'
$ws.Cells.Item(285, 1).Value = '13:35:25'
$ws.Cells.Item(285, 2).Value = '15:17'
$ws.Cells.Item(285, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(285, 4).Value = 102
$ws.Cells.Item(285, 5).Value = 'LP1912'
$ws.Cells.Item(286, 1).Value = '13:35:25'
$ws.Cells.Item(286, 2).Value = '15:24'
$ws.Cells.Item(286, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(286, 4).Value = 109
$ws.Cells.Item(286, 5).Value = 'LP1912'
$ws.Cells.Item(287, 1).Value = '14:23:38'
$ws.Cells.Item(287, 2).Value = '15:24'
$ws.Cells.Item(287, 3).Value = '215C_EL PATO'
$ws.Cells.Item(287, 4).Value = 61
$ws.Cells.Item(287, 5).Value = 'LP1912'
$ws.Cells.Item(288, 1).Value = '13:35:25'
$ws.Cells.Item(288, 2).Value = '15:25'
$ws.Cells.Item(288, 3).Value = '215C_EL PATO'
$ws.Cells.Item(288, 4).Value = 110
$ws.Cells.Item(288, 5).Value = 'LP1912'
$ws.Cells.Item(289, 1).Value = '13:54:15'
$ws.Cells.Item(289, 2).Value = '15:25'
$ws.Cells.Item(289, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(289, 4).Value = 91
$ws.Cells.Item(289, 5).Value = 'LP1912'
$ws.Cells.Item(290, 1).Value = '13:54:15'
$ws.Cells.Item(290, 2).Value = '15:25'
$ws.Cells.Item(290, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(290, 4).Value = 91
$ws.Cells.Item(290, 5).Value = 'LP1912'
$ws.Cells.Item(291, 1).Value = '13:54:15'
$ws.Cells.Item(291, 2).Value = '15:36'
$ws.Cells.Item(291, 3).Value = '17X38_ROMERO'
$ws.Cells.Item(291, 4).Value = 102
$ws.Cells.Item(291, 5).Value = 'LP1912'
$ws.Cells.Item(292, 1).Value = '14:23:38'
$ws.Cells.Item(292, 2).Value = '15:44'
$ws.Cells.Item(292, 3).Value = '14_ABASTO'
$ws.Cells.Item(292, 4).Value = 81
$ws.Cells.Item(292, 5).Value = 'LP1912'
$ws.Cells.Item(293, 1).Value = '14:23:38'
$ws.Cells.Item(293, 2).Value = '15:45'
$ws.Cells.Item(293, 3).Value = '215C_LA PLATA'
$ws.Cells.Item(293, 4).Value = 82
$ws.Cells.Item(293, 5).Value = 'LP1912'
$ws.Cells.Item(294, 1).Value = '14:23:38'
$ws.Cells.Item(294, 2).Value = '15:55'
$ws.Cells.Item(294, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(294, 4).Value = 92
$ws.Cells.Item(294, 5).Value = 'LP1912'
$ws.Cells.Item(295, 1).Value = '14:23:38'
$ws.Cells.Item(295, 2).Value = '16:01'
$ws.Cells.Item(295, 3).Value = '15_ABASTO'
$ws.Cells.Item(295, 4).Value = 98
$ws.Cells.Item(295, 5).Value = 'LP1912'
$ws.Cells.Item(296, 1).Value = '14:23:38'
$ws.Cells.Item(296, 2).Value = '16:20'
$ws.Cells.Item(296, 3).Value = '10_OLMOS'
$ws.Cells.Item(296, 4).Value = 117
$ws.Cells.Item(296, 5).Value = 'LP1912'
